$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.813.45"
$ws.Range("E2").Value = "  -0.13%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.639.87"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.19%  "

# Row 5 - BNB (purely-numeric text price -> force text with a quote prefix)
$ws.Range("D5").Formula = "'215.93"
$ws.Range("E5").Value = "  +0.46%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.58%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.19%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.26%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.96%  "

# Row 10 - Solana
$ws.Range("D10").Formula = "'19.69"
$ws.Range("E10").Value = "  -0.92%  "

# Row 11 - TRON
$ws.Range("D11").Formula = "'0.0792"
$ws.Range("E11").Value = "  +1.40%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  +0.24%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "1.865.61"
$ws.Range("E13").Value = "  +0.31%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.641.43"
$ws.Range("E14").Value = "  +0.39%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.45%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.12%  "

# Row 17 - Litecoin
$ws.Range("D17").Formula = "'63.10"
$ws.Range("E17").Value = "  +0.12%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.843.99"
$ws.Range("E18").Value = "  -0.09%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.17%  "

# Row 20 - Uniswap
$ws.Range("D20").Formula = "'4.48"
$ws.Range("E20").Value = "  +2.14%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Formula = "'192.54"
$ws.Range("E21").Value = "  -0.45%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  +0.67%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  +2.35%  "

# Row 24 - Toncoin
$ws.Range("D24").Formula = "'1.85"
$ws.Range("E24").Value = "  +5.18%  "

# Row 25 - BinanceUSD
$ws.Range("E25").Value = "  -0.14%  "

# Row 26 - Monero
$ws.Range("D26").Formula = "'142.30"
$ws.Range("E26").Value = "  +2.18%  "

# Row 27 - Stellar
$ws.Range("D27").Formula = "'0.123"
$ws.Range("E27").Value = "  +1.29%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +1.87%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.33%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.26%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +0.93%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -0.19%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +0.54%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -0.07%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Formula = "'0.907"
$ws.Range("E36").Value = "  +0.78%  "

# Row 37 - Maker
$ws.Range("D37").Value = "1.133.20"

# Row 38 - MXToken
$ws.Range("E38").Value = "  -1.76%  "

# Row 39 - ImmutableX
$ws.Range("D39").Formula = "'0.545"
$ws.Range("E39").Value = "  -0.71%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +0.29%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.05%  "

# Row 42 - FraxShare
$ws.Range("E42").Value = "  +1.36%  "

# Row 43 - Quant
$ws.Range("D43").Formula = "'100.75"
$ws.Range("E43").Value = "  +1.13%  "

# Row 44 - TrustWalletToken
$ws.Range("E44").Value = "  +0.69%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.775.32"
$ws.Range("E45").Value = "  +0.07%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  +3.07%  "

# Row 47 - Aave
$ws.Range("D47").Formula = "'55.38"
$ws.Range("E47").Value = "  -0.07%  "

# Rows 49/50 swap - RenderToken <-> Cronos
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Formula = "'0.0502"
$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Formula = "'1.43"
$ws.Range("E50").Value = "  +4.65%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  +1.72%  "
